$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nbsp = [char]0x00A0

# Row 223
$ws.Range("A3:C3").Copy($ws.Range("A223:C223"))
$ws.Rows.Item(223).RowHeight = 18
$ws.Cells.Item(223, 1).Value = "Tuesday, April 01, 2025"
$ws.Cells.Item(223, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(223, 3)
$cC.Value = "Job Openings and Labor Turnover Survey" + $nbsp + "for February 2025"
$titleLen = "Job Openings and Labor Turnover Survey".Length
$suffixText = $nbsp + "for February 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 224
$ws.Range("A6:C6").Copy($ws.Range("A224:C224"))
$ws.Rows.Item(224).RowHeight = 18
$ws.Cells.Item(224, 1).Value = "Wednesday, April 02, 2025"
$ws.Cells.Item(224, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(224, 3)
$cC.Value = "Occupational Employment and Wages" + $nbsp + "for Annual 2024"
$titleLen = "Occupational Employment and Wages".Length
$suffixText = $nbsp + "for Annual 2024"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 225
$ws.Range("A3:C3").Copy($ws.Range("A225:C225"))
$ws.Rows.Item(225).RowHeight = 18
$ws.Cells.Item(225, 1).Value = "Friday, April 04, 2025"
$ws.Cells.Item(225, 2).Value = 0.35416666666666669
$cC = $ws.Cells.Item(225, 3)
$cC.Value = "Employment Situation" + $nbsp + "for March 2025"
$titleLen = "Employment Situation".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 226
$ws.Range("A4:C4").Copy($ws.Range("A226:C226"))
$ws.Rows.Item(226).RowHeight = 27
$ws.Cells.Item(226, 1).Value = "Wednesday, April 09, 2025"
$ws.Cells.Item(226, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(226, 3)
$cC.Value = "Metropolitan Area Employment and Unemployment (Monthly)" + $nbsp + "for February 2025"
$titleLen = "Metropolitan Area Employment and Unemployment (Monthly)".Length
$suffixText = $nbsp + "for February 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 227
$ws.Range("A3:C3").Copy($ws.Range("A227:C227"))
$ws.Rows.Item(227).RowHeight = 18
$ws.Cells.Item(227, 1).Value = "Thursday, April 10, 2025"
$ws.Cells.Item(227, 2).Value = 0.35416666666666669
$cC = $ws.Cells.Item(227, 3)
$cC.Value = "Consumer Price Index" + $nbsp + "for March 2025"
$titleLen = "Consumer Price Index".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 228
$ws.Range("A52:C52").Copy($ws.Range("A228:C228"))
$ws.Cells.Item(228, 1).Value = "Thursday, April 10, 2025"
$ws.Cells.Item(228, 2).Value = 0.35416666666666669
$cC = $ws.Cells.Item(228, 3)
$cC.Value = "Real Earnings" + $nbsp + "for March 2025"
$titleLen = "Real Earnings".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 229
$ws.Range("A7:C7").Copy($ws.Range("A229:C229"))
$ws.Cells.Item(229, 1).Value = "Friday, April 11, 2025"
$ws.Cells.Item(229, 2).Value = 0.35416666666666669
$cC = $ws.Cells.Item(229, 3)
$cC.Value = "Producer Price Index" + $nbsp + "for March 2025"
$titleLen = "Producer Price Index".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 230
$ws.Range("A6:C6").Copy($ws.Range("A230:C230"))
$ws.Rows.Item(230).RowHeight = 18
$ws.Cells.Item(230, 1).Value = "Tuesday, April 15, 2025"
$ws.Cells.Item(230, 2).Value = 0.35416666666666669
$cC = $ws.Cells.Item(230, 3)
$cC.Value = "U.S. Import and Export Price Indexes" + $nbsp + "for March 2025"
$titleLen = "U.S. Import and Export Price Indexes".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 231
$ws.Range("A3:C3").Copy($ws.Range("A231:C231"))
$ws.Rows.Item(231).RowHeight = 18
$ws.Cells.Item(231, 1).Value = "Wednesday, April 16, 2025"
$ws.Cells.Item(231, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(231, 3)
$cC.Value = "State Job Openings and Labor Turnover" + $nbsp + "for February 2025"
$titleLen = "State Job Openings and Labor Turnover".Length
$suffixText = $nbsp + "for February 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 232
$ws.Range("A6:C6").Copy($ws.Range("A232:C232"))
$ws.Rows.Item(232).RowHeight = 18
$ws.Cells.Item(232, 1).Value = "Wednesday, April 16, 2025"
$ws.Cells.Item(232, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(232, 3)
$cC.Value = "Usual Weekly Earnings of Wage and Salary Workers" + $nbsp + "for First Quarter 2025"
$titleLen = "Usual Weekly Earnings of Wage and Salary Workers".Length
$suffixText = $nbsp + "for First Quarter 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 233
$ws.Range("A13:C13").Copy($ws.Range("A233:C233"))
$ws.Rows.Item(233).RowHeight = 27
$ws.Cells.Item(233, 1).Value = "Friday, April 18, 2025"
$ws.Cells.Item(233, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(233, 3)
$cC.Value = "State Employment and Unemployment (Monthly)" + $nbsp + "for March 2025"
$titleLen = "State Employment and Unemployment (Monthly)".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 234
$ws.Range("A4:C4").Copy($ws.Range("A234:C234"))
$ws.Rows.Item(234).RowHeight = 27
$ws.Cells.Item(234, 1).Value = "Tuesday, April 22, 2025"
$ws.Cells.Item(234, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(234, 3)
$cC.Value = "College Enrollment and Work Activity of High School Graduates" + $nbsp + "for Annual 2024"
$titleLen = "College Enrollment and Work Activity of High School Graduates".Length
$suffixText = $nbsp + "for Annual 2024"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 235
$ws.Range("A3:C3").Copy($ws.Range("A235:C235"))
$ws.Rows.Item(235).RowHeight = 18
$ws.Cells.Item(235, 1).Value = "Wednesday, April 23, 2025"
$ws.Cells.Item(235, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(235, 3)
$cC.Value = "Employment Characteristics of Families" + $nbsp + "for Annual 2024"
$titleLen = "Employment Characteristics of Families".Length
$suffixText = $nbsp + "for Annual 2024"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 236
$ws.Range("A6:C6").Copy($ws.Range("A236:C236"))
$ws.Rows.Item(236).RowHeight = 18
$ws.Cells.Item(236, 1).Value = "Tuesday, April 29, 2025"
$ws.Cells.Item(236, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(236, 3)
$cC.Value = "Job Openings and Labor Turnover Survey" + $nbsp + "for March 2025"
$titleLen = "Job Openings and Labor Turnover Survey".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 237
$ws.Range("A13:C13").Copy($ws.Range("A237:C237"))
$ws.Rows.Item(237).RowHeight = 27
$ws.Cells.Item(237, 1).Value = "Tuesday, April 29, 2025"
$ws.Cells.Item(237, 2).Value = 0.41666666666666669
$cC = $ws.Cells.Item(237, 3)
$cC.Value = "Metropolitan Area Employment and Unemployment (Monthly)" + $nbsp + "for March 2025"
$titleLen = "Metropolitan Area Employment and Unemployment (Monthly)".Length
$suffixText = $nbsp + "for March 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

# Row 238
$ws.Range("A6:C6").Copy($ws.Range("A238:C238"))
$ws.Rows.Item(238).RowHeight = 18
$ws.Cells.Item(238, 1).Value = "Wednesday, April 30, 2025"
$ws.Cells.Item(238, 2).Value = 0.35416666666666669
$cC = $ws.Cells.Item(238, 3)
$cC.Value = "Employment Cost Index" + $nbsp + "for First Quarter 2025"
$titleLen = "Employment Cost Index".Length
$suffixText = $nbsp + "for First Quarter 2025"
$runChars = $cC.Characters($titleLen + 1, $suffixText.Length)
$runChars.Font.Bold = $false
$runChars.Font.Size = 7
$runChars.Font.Name = "Tahoma"
$runChars.Font.Color = 3355443

$null = $ws.Range("E223").Select()
